# Quarterly database update for Overview sheet.
# - Rolls the 10 visible quarter columns (E:N) one quarter forward:
#     the oldest quarter (E) is dropped, every column shifts left by one,
#     and the newest quarter ("فصل اول منتهی به 1401/11") is appended in N.
# - Applies the same roll-forward to the quarter-header rows (8 and 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Quarter header labels (rows 8 and 24, columns E:N) ----------------
$quarterHeaders = @(
    "فصل اول منتهی به 1399/08",
    "فصل چهارم منتهی به 1399/12",
    "فصل سوم منتهی به 1400/05",
    "فصل دوم منتهی به 1400/06",
    "فصل چهارم منتهی به 1400/08",
    "فصل اول منتهی به 1400/11",
    "فصل دوم منتهی به 1401/02",
    "فصل سوم منتهی به 1401/05",
    "فصل چهارم منتهی به 1401/08",
    "فصل اول منتهی به 1401/11"
)

foreach ($headerRow in @(8, 24)) {
    for ($i = 0; $i -lt $quarterHeaders.Length; $i++) {
        $col = 5 + $i   # column E = 5
        $ws.Cells.Item($headerRow, $col).Value = $quarterHeaders[$i]
    }
}

# ---- 2. Data rows: shift quarter values left by one, append new quarter --
# Each entry is the *final* (post-roll-forward) row of values for columns
# E..N. "-" marks the placeholder/no-data cells (shared "-" text already
# used throughout the sheet).
$dataRows = @{
    10 = @(0, 0, "-", 0, "-", 0, 0, 0, 0, 0)
    11 = @(0, 0, "-", 0, "-", 0, 0, 0, 0, 0)
    12 = @(406, 12139, "-", 9733, "-", 0, 0, 0, 0, 0)
    13 = @(398, 785, "-", 159, "-", 0, 1162, -1162, 0, 0)
    14 = @(0, 0, "-", 0, "-", 0, 0, 0, 0, 0)
    15 = @(0, 0, "-", 0, "-", 0, 0, 1981, -1981, 0)
    16 = @(326, 1589, "-", 2124, "-", 963, 1873, 887, 1999, 2143)
    17 = @(14458, 28516, "-", 31735, "-", 24900, 54931, 39769, 44049, 44014)
    18 = @(0, 0, "-", 0, "-", 0, 0, 0, 0, 0)
    19 = @(15752, 53352, "-", 76437, "-", 69861, 84804, 150053, 22766, 105560)
    20 = @(31340, 96381, 0, 120188, 0, 95724, 142770, 191528, 66833, 151717)
    26 = @(187, 132, "-", 201, "-", 257, 257, 257, 255, 255)
    27 = @(125, 178, "-", 130, "-", 89, 89, 89, 92, 92)
}

foreach ($rowNum in $dataRows.Keys) {
    $values = $dataRows[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i   # column E = 5 .. N = 14
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}

Write-Output "quarterly roll-forward applied"
